# The document had a number of paragraphs that ended with one extra,
# redundant tab character (runs of <w:tab/> elements). This script removes
# exactly one tab from each such run, matching the "optimized" document.
#
# Because several paragraphs are not unique by text (the tab-only ones),
# we operate by paragraph index and, within each paragraph, remove one
# specific occurrence of the tab character (counted left-to-right across
# the whole paragraph).

$d = $word.ActiveDocument

function Remove-TabOccurrence($paraIndex, $occurrence) {
    $para = $d.Paragraphs.Item($paraIndex)
    $rng = $para.Range
    $text = $rng.Text
    $seen = 0
    for ($i = 0; $i -lt $text.Length; $i++) {
        if ($text.Substring($i, 1) -eq "`t") {
            $seen = $seen + 1
            if ($seen -eq $occurrence) {
                $pos = $rng.Start + $i
                $tabRange = $d.Range($pos, $pos + 1)
                $tabRange.Delete()
                return $true
            }
        }
    }
    return $false
}

# "Название предприятия: " tab run -> drop the trailing (10th) tab
Remove-TabOccurrence 5 10 | Out-Null
# Blank continuation line below it -> drop the trailing (13th) tab
Remove-TabOccurrence 6 13 | Out-Null

# "Почтовый адрес: " tab run -> drop the trailing (11th) tab
Remove-TabOccurrence 9 11 | Out-Null
# Blank continuation line below it -> drop the trailing (13th) tab
Remove-TabOccurrence 10 13 | Out-Null

# "Телефон, телефон/факс, e-mail: " tab run -> drop the trailing (8th) tab
Remove-TabOccurrence 13 8 | Out-Null
# Blank continuation line below it -> drop the trailing (13th) tab
Remove-TabOccurrence 14 13 | Out-Null

# "Ф.И.О. директора: " tab run -> drop the trailing (10th) tab
Remove-TabOccurrence 17 10 | Out-Null
# Blank continuation line below it -> drop the trailing (13th) tab
Remove-TabOccurrence 18 13 | Out-Null

# "Председатель" line: tab, "Председатель", tab, tab -> drop the tab
# immediately following the word (2nd tab overall in the paragraph)
Remove-TabOccurrence 30 2 | Out-Null

# Signature line: 6 leading tabs before "(подпись)" -> drop the first one
Remove-TabOccurrence 31 1 | Out-Null
